$wb2 = $excel.ActiveWorkbook
$ws = $wb2.Worksheets.Item(1)

# Rows 2 and 3 in this sighting export describe two different species recorded
# at the same locality/date/observer. The edit swaps which row holds which
# species' data (Id, Taxonsorteringsordning, Rödlistade, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Enhet, Kön, Aktivitet, Bestämningsmetod) while
# every shared column (locality, coordinates, date, observer, substrate, ...)
# is identical between the two rows and is left untouched.
$cols = @("A","B","D","E","F","G","H","J","L","M","AF")

foreach ($col in $cols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value()
    $v3 = $ws.Range($addr3).Value()
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# After the swap, row 2 needs an (empty) "Kön" cell and row 3 needs (empty)
# "Enhet"/"Bestämningsmetod" cells to exist (they were present-but-blank on
# the row the data came from). A bare "" assignment on a previously absent
# cell doesn't materialise it, so briefly mark the cell as text (leading
# apostrophe) and then restore the default "Normal" style.
$blankCells = @("L2", "J3", "AF3")
foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}
